$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old C11 text ("...tauschseite und anfragenseite(noch ohne datenbankanbindung)") is being
# split: the "anfragenseite" part becomes its own, newly logged row (12), while C11 keeps the
# shorter "accountmanage...tauschseite" text. Clear C11 first so the freed shared-string slot
# gets reused by the new row's text, matching how the shared-string table is laid out upstream.
$ws.Range("C11").ClearContents()

# New row 12: continue the time log with the entry about the anfragenseite work.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats - reuse A11's date format
$ws.Range("A12").Value = 42732
$ws.Range("B12").Value = 6
$ws.Range("C12").Value = "anfragenseite komplettes layout und controllerfunktionen,detailseite gefüllt"

# Now give C11 its shortened text.
$ws.Range("C11").Value = "accountmanage seit mit passwort ändern,reservierungen löschbar, details,tauschseite"

# Match the selection reflected in the saved workbook.
$ws.Range("G12").Select()
